$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.760000000000001
$ws.Range("A4").Value = -21.352
$ws.Range("B4").Value = 5.059
$ws.Range("D4").Value = -7.964000000000001
$ws.Range("A6").Value = -21.527
$ws.Range("C6").Value = -12.775
$ws.Range("A7").Value = -20.451
$ws.Range("C7").Value = -12.154
$ws.Range("A8").Value = -21.204
$ws.Range("B8").Value = 5.958
$ws.Range("C8").Value = -13.282
$ws.Range("B9").Value = 5.52
$ws.Range("C10").Value = -12.965
$ws.Range("B12").Value = 5.424
$ws.Range("C13").Value = -12.784
$ws.Range("D13").Value = -7.645000000000001
$ws.Range("C14").Value = -13.231
$ws.Range("A16").Value = -21.779
$ws.Range("C16").Value = -12.598
$ws.Range("D16").Value = -7.343999999999999
$ws.Range("B17").Value = 5.534000000000001
$ws.Range("B18").Value = 5.51
$ws.Range("B19").Value = 7.026999999999999
$ws.Range("A20").Value = -21.711
$ws.Range("B20").Value = 6.417999999999999
$ws.Range("D20").Value = -7.494999999999999
$ws.Range("A21").Value = -20.186
$ws.Range("D25").Value = -8.107000000000001
$ws.Range("B26").Value = 5.69
$ws.Range("A28").Value = -21.555
$ws.Range("A29").Value = -21.128
$ws.Range("A30").Value = -22.172
$ws.Range("C30").Value = -12.657
$ws.Range("B31").Value = 6.06
$ws.Range("A32").Value = -21.563
$ws.Range("D34").Value = -7.715999999999999
$ws.Range("C37").Value = -13.013
$ws.Range("B39").Value = 7.043000000000001
$ws.Range("D39").Value = -7.755
$ws.Range("A40").Value = -21.587
$ws.Range("B40").Value = 5.786
$ws.Range("C40").Value = -11.26
$ws.Range("B41").Value = 7.362
$ws.Range("B42").Value = 6.286
$ws.Range("B43").Value = 5.687
$ws.Range("C44").Value = -12.179
$ws.Range("A46").Value = -21.606
$ws.Range("B47").Value = 5.833
$ws.Range("B48").Value = 5.441
$ws.Range("A51").Value = -21.585
$ws.Range("D51").Value = -8.32
$ws.Range("A52").Value = -21.463
$ws.Range("B54").Value = 5.467000000000001
$ws.Range("A57").Value = -21.806
$ws.Range("A59").Value = -21.678
$ws.Range("D59").Value = -7.688999999999998
$ws.Range("D61").Value = -8.019000000000002
$ws.Range("A62").Value = -21.502
$ws.Range("B62").Value = 5.951000000000001
$ws.Range("D62").Value = -7.915000000000001
$ws.Range("B63").Value = 5.334000000000001
$ws.Range("B64").Value = 5.585000000000001
$ws.Range("D64").Value = -7.483
$ws.Range("A66").Value = -21.272
$ws.Range("D69").Value = -7.938999999999998
$ws.Range("C70").Value = -11.676
$ws.Range("A73").Value = -21.454
$ws.Range("A74").Value = -20.935
$ws.Range("B76").Value = 6.778
$ws.Range("A77").Value = -20.952
$ws.Range("D78").Value = -8.478
$ws.Range("B81").Value = 5.982000000000001
$ws.Range("D83").Value = -8.286999999999999
$ws.Range("B84").Value = 6.161
$ws.Range("B89").Value = 5.403
$ws.Range("C89").Value = -13.259
$ws.Range("C91").Value = -12.537
$ws.Range("A92").Value = -21.492
$ws.Range("D92").Value = -7.222
$ws.Range("C93").Value = -11.625
$ws.Range("B94").Value = 5.912999999999999
$ws.Range("C98").Value = -12.089
$ws.Range("D98").Value = -7.737
$ws.Range("A100").Value = -21.286
$ws.Range("D100").Value = -7.806
